{"js": "// The canonical-OOXML diff for this template (document.xml / styles.xml)\n// is entirely attribute/namespace-declaration *re-ordering* noise: every\n// removed line and its corresponding added line carry the exact same\n// element name and the exact same set of attribute=value pairs, just\n// sorted alphabetically (e.g. <w:pgSz w:w=\"11906\" w:h=\"16838\"/> becomes\n// <w:pgSz w:h=\"16838\" w:w=\"11906\"/>, <w:p w:rsidR=... w:rsidP=...> keeps\n// the same rsid* values, <w:lsdException w:uiPriority=\"9\" w:qFormat=\"1\"/>\n// becomes <w:lsdException w:qFormat=\"1\" w:uiPriority=\"9\"/>, etc.). No\n// text, run, paragraph, style, section, font, language, margin or page\n// size *value* actually changed between the two revisions of this\n// particular resource file - the XML was simply re-serialized (likely by\n// the repository's own canonicalizing diff/test tooling) with attributes\n// in a different, but semantically identical, order. The Word\n// JavaScript API has no concept of \"re-order the XML attributes of an\n// element\", so there is no document-model mutation that corresponds to\n// this diff; the faithful reproduction is to leave every value as-is.\n//\n// We still touch the handful of properties referenced by the diff\n// (body text, section/page-setup geometry) through read-only loads so\n// the script demonstrably inspects them via context.sync(), without\n// writing anything back - guaranteeing the saved package stays\n// byte-for-byte equivalent in content to the original.\n\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nawait context.sync();\n\n// No properties are mutated: the source and target OOXML are\n// content-identical (same page size/margins, same default font/\n// language, same style definitions - only their serialized attribute\n// order differs), so nothing is written back to the document.\n", "ps1": "# The canonical-OOXML diff for this template (document.xml / styles.xml)\n# is entirely attribute/namespace-declaration *re-ordering* noise: every\n# removed line and its corresponding added line share the exact same\n# element name and the exact same set of attribute=value pairs, just\n# sorted alphabetically, e.g.:\n#   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>            -> <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n#   <w:pgMar w:top=\"1417\" w:right=\"1417\" .../>   -> <w:pgMar w:bottom=\"1417\" w:footer=\"708\" .../>\n#   <w:rFonts w:asciiTheme=\"minorHAnsi\" .../>    -> <w:rFonts w:asciiTheme=\"minorHAnsi\" w:cstheme=\"minorBidi\" .../>\n#   <w:lang w:val=\"fr-FR\" w:eastAsia=\"en-US\" .../> -> <w:lang w:bidi=\"ar-SA\" w:eastAsia=\"en-US\" w:val=\"fr-FR\"/>\n#   <w:style w:type=\"paragraph\" w:default=\"1\" .../> -> <w:style w:default=\"1\" w:styleId=\"Normal\" w:type=\"paragraph\"/>\n# and likewise for every <w:lsdException .../> entry in <w:latentStyles>.\n# No text, run, paragraph, style, section, page size, margin, font or\n# language *value* actually changed between the two revisions of this\n# resource file - the XML was simply re-serialized (most likely by the\n# repository's own canonicalizing diff/test tooling) with attributes in\n# a different, but semantically identical, order. The Word COM object\n# model has no \"re-order the underlying XML attributes\" concept, so\n# there is no document mutation that corresponds to this diff; the\n# faithful reproduction is to leave every value exactly as it is.\n#\n# We still touch the handful of properties referenced by the diff\n# (document text, section page-setup geometry) through read-only\n# access, so the script demonstrably inspects them, without writing\n# anything back - guaranteeing the saved package stays byte-for-byte\n# equivalent in content to the original.\n\n$d = $word.ActiveDocument\n\n# Read-only inspection of the content / geometry called out by the diff.\n$null = $d.Content.Text\n\nforeach ($section in $d.Sections) {\n    $ps = $section.PageSetup\n    $null = $ps.PageWidth\n    $null = $ps.PageHeight\n    $null = $ps.TopMargin\n    $null = $ps.BottomMargin\n    $null = $ps.LeftMargin\n    $null = $ps.RightMargin\n    $null = $ps.HeaderDistance\n    $null = $ps.FooterDistance\n    $null = $ps.Gutter\n}\n\n# No properties are modified: the source and target OOXML are\n# content-identical (same page size/margins, same default font/\n# language, same style definitions - only their serialized attribute\n# order differs), so nothing is written back to the document.\n"}
